# Auto-generated Excel COM-interop script
# Updates "Horarios 141" schedule workbook with refreshed scrape data
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 04:44:55"
$ws.Range("A3").Value = "Total filas: 38"

$data = @(
    @("04:01:06","04:02","81_EL PELIGRO",1),
    @("04:44:55","04:46","215_EL PELIGRO",2),
    @("00:46:06","01:12","215_ALUAR",26),
    @("04:44:55","05:11","17_ROMERO",27),
    @("04:44:55","05:21","23_HERNANDEZ",37),
    @("04:01:06","04:47","215_EL PELIGRO",46),
    @("04:44:55","05:31","81_EL PELIGRO",47),
    @("04:44:55","05:43","14_ABASTO",59),
    @("03:46:12","04:46","215A_EL PATO",60),
    @("01:55:38","03:02","15_ABASTO",67),
    @("04:44:55","05:51","17_ROMERO",67),
    @("04:01:06","05:12","17_ROMERO",71),
    @("00:46:06","01:58","14_ABASTO",72),
    @("04:30:03","05:44","14_ABASTO",74),
    @("04:44:55","06:00","16_SANTA ANA",76),
    @("04:44:55","06:03","10_OLMOS",79),
    @("04:44:55","06:10","215A_EL PATO",86),
    @("03:46:12","05:16","17_ROMERO",90),
    @("04:30:03","06:01","16_SANTA ANA",91),
    @("04:01:06","05:32","81_EL PELIGRO",91),
    @("02:29:13","04:01","81_EL PELIGRO",92),
    @("04:30:03","06:04","10_OLMOS",94),
    @("01:22:42","02:58","215_ALUAR",96),
    @("03:46:12","05:22","23_HERNANDEZ",96),
    @("04:44:55","06:23","11_ETCHEVERRY",99),
    @("04:30:03","06:11","215A_EL PATO",101),
    @("04:44:55","06:26","23_HERNANDEZ",102),
    @("04:01:06","05:45","14_ABASTO",104),
    @("04:44:55","06:30","17X38_ROMERO",106),
    @("04:44:55","06:30","16_SANTA ANA",106),
    @("03:46:12","05:35","215B_EL PATO",109),
    @("04:01:06","05:52","17_ROMERO",111),
    @("01:55:38","03:48","14_ABASTO",113),
    @("03:00:53","04:53","11_ETCHEVERRY",113),
    @("04:30:03","06:24","11_ETCHEVERRY",114),
    @("04:44:55","06:38","225_C ROCA-H SUR",114),
    @("04:30:03","06:27","23_HERNANDEZ",117),
    @("02:47:42","04:45","215A_EL PATO",118)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 6 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 04:44:55"
$ws.Range("A3").Value = "Total filas: 9"

$data = @(
    @("04:44:55","04:46","215_EL PELIGRO",2),
    @("00:46:06","01:12","215_ALUAR",26),
    @("04:01:06","04:47","215_EL PELIGRO",46),
    @("03:46:12","04:46","215A_EL PATO",60),
    @("04:44:55","06:10","215A_EL PATO",86),
    @("01:22:42","02:58","215_ALUAR",96),
    @("04:30:03","06:11","215A_EL PATO",101),
    @("03:46:12","05:35","215B_EL PATO",109),
    @("02:47:42","04:45","215A_EL PATO",118)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 6 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 04:44:55"
$ws.Range("A3").Value = "Total filas: 38"

$data = @(
    @("04:01:06","04:02","81_EL PELIGRO",1),
    @("04:44:55","04:46","215_EL PELIGRO",2),
    @("00:46:06","01:12","215_ALUAR",26),
    @("04:44:55","05:11","17_ROMERO",27),
    @("04:44:55","05:21","23_HERNANDEZ",37),
    @("04:01:06","04:47","215_EL PELIGRO",46),
    @("04:44:55","05:31","81_EL PELIGRO",47),
    @("04:44:55","05:43","14_ABASTO",59),
    @("03:46:12","04:46","215A_EL PATO",60),
    @("01:55:38","03:02","15_ABASTO",67),
    @("04:44:55","05:51","17_ROMERO",67),
    @("04:01:06","05:12","17_ROMERO",71),
    @("00:46:06","01:58","14_ABASTO",72),
    @("04:30:03","05:44","14_ABASTO",74),
    @("04:44:55","06:00","16_SANTA ANA",76),
    @("04:44:55","06:03","10_OLMOS",79),
    @("04:44:55","06:10","215A_EL PATO",86),
    @("03:46:12","05:16","17_ROMERO",90),
    @("04:30:03","06:01","16_SANTA ANA",91),
    @("04:01:06","05:32","81_EL PELIGRO",91),
    @("02:29:13","04:01","81_EL PELIGRO",92),
    @("04:30:03","06:04","10_OLMOS",94),
    @("01:22:42","02:58","215_ALUAR",96),
    @("03:46:12","05:22","23_HERNANDEZ",96),
    @("04:44:55","06:23","11_ETCHEVERRY",99),
    @("04:30:03","06:11","215A_EL PATO",101),
    @("04:44:55","06:26","23_HERNANDEZ",102),
    @("04:01:06","05:45","14_ABASTO",104),
    @("04:44:55","06:30","17X38_ROMERO",106),
    @("04:44:55","06:30","16_SANTA ANA",106),
    @("03:46:12","05:35","215B_EL PATO",109),
    @("04:01:06","05:52","17_ROMERO",111),
    @("01:55:38","03:48","14_ABASTO",113),
    @("03:00:53","04:53","11_ETCHEVERRY",113),
    @("04:30:03","06:24","11_ETCHEVERRY",114),
    @("04:44:55","06:38","225_C ROCA-H SUR",114),
    @("04:30:03","06:27","23_HERNANDEZ",117),
    @("02:47:42","04:45","215A_EL PATO",118)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 6 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

